$wb = $excel.ActiveWorkbook

# Add a new worksheet "InvalidLogin" after the last existing sheet (ValidLogin)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "InvalidLogin"

# Populate the new sheet with the invalid-login test data (data driven test)
$ws3.Range("A1").Value = "username"
$ws3.Range("B1").Value = "password"
$ws3.Range("A2").Value = "abc"
$ws3.Range("B2").Value = "xyz"

# Match authored selection state: InvalidLogin tab active, B2 selected
$ws3.Range("B2").Select() | Out-Null

# ValidLogin sheet is no longer the active/selected tab; its selection
# moves to the A1:B1 header row
$ws2 = $wb.Worksheets.Item("ValidLogin")
$ws2.Range("A1:B1").Select() | Out-Null

# InvalidLogin becomes the active sheet/tab
$ws3.Activate() | Out-Null
